# Commit: "Implemented Command of executors"
#
# The deck originally contains a single slide (the "UI Design" workspace
# diagram: Registration / Browser / Workspaces / Scheduler, etc.). The
# canonical edit duplicates that slide and places the copy in front of the
# original, so the presentation ends up with two (structurally identical)
# slides: slide 256 keeps its original relationship/content and slide 257
# is the newly duplicated slide that is inserted right after it, both
# carrying the same "UI Design" diagram content.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Duplicate slide 1 -- PowerPoint inserts the duplicate immediately after
# the source slide, which reproduces the target slide order (original
# slide id 256 first, new duplicate id 257 second).
$s.Duplicate() | Out-Null
